$d = $word.ActiveDocument

$pairs = @(
    @("357÷7=51, 0", "734÷5=146, 4"),
    @("305÷2=152, 1", "318÷2=159, 0"),
    @("293÷2=146, 1", "247÷8=30, 7"),
    @("847÷5=169, 2", "431÷5=86, 1"),
    @("356÷3=118, 2", "345÷5=69, 0"),
    @("645÷3=215, 0", "356÷5=71, 1"),
    @("385÷7=55, 0", "352÷2=176, 0"),
    @("661÷9=73, 4", "652÷4=163, 0"),
    @("573÷8=71, 5", "853÷7=121, 6"),
    @("760÷4=190, 0", "854÷4=213, 2"),
    @("717÷9=79, 6", "301÷9=33, 4"),
    @("704÷6=117, 2", "224÷2=112, 0"),
    @("989÷4=247, 1", "374÷7=53, 3"),
    @("329÷3=109, 2", "923÷6=153, 5"),
    @("798÷6=133, 0", "280÷8=35, 0"),
    @("565÷6=94, 1", "322÷4=80, 2"),
    @("304÷7=43, 3", "870÷5=174, 0"),
    @("142÷2=71, 0", "620÷3=206, 2"),
    @("364÷4=91, 0", "667÷4=166, 3"),
    @("955÷7=136, 3", "374÷9=41, 5"),
    @("433÷3=144, 1", "514÷2=257, 0"),
    @("855÷6=142, 3", "989÷8=123, 5"),
    @("971÷8=121, 3", "372÷7=53, 1"),
    @("942÷5=188, 2", "621÷3=207, 0"),
    @("315÷8=39, 3", "920÷3=306, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"
